# Append a new job listing at row 13 of the "ランサーズ" sheet, pushing the
# previous rows 13-17 down to 14-18, refresh the "取得日時" (fetched-at)
# timestamp for every data row, widen column D, and rebuild the F-column
# hyperlinks so they keep pointing at the right row after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-30 01:21:17"

# --- 1. Shift rows 13..17 down to 14..18 (bottom-up so we never overwrite
#        a source row before it has been read). Columns A-G only; column H
#        is blank for all of these rows both before and after the edit. ---
for ($r = 17; $r -ge 13; $r--) {
    $dst = $r + 1
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($dst, $c).Value2 = $ws.Cells.Item($r, $c).Value2
    }
}

# --- 2. Write the brand-new row 13 ---
$ws.Cells.Item(13, 1).Value2 = $newTimestamp
$ws.Cells.Item(13, 2).Value2 = "【フルリモート】SESエンジニア募集|スキルに応じて30〜40万円/月|複数案件あり・継続前提"
$ws.Cells.Item(13, 3).Value2 = "システム開発"
$ws.Cells.Item(13, 4).Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(13, 5).Value2 = "期限情報なし"
$ws.Cells.Item(13, 6).Value2 = "https://www.lancers.jp/work/detail/5417644"
$ws.Cells.Item(13, 7).Value2 = 25

# --- 3. Row 18 (old row 17 content) also picked up a wording change in
#        column D beyond the plain shift. ---
$ws.Cells.Item(18, 4).Value2 = "10,000 円 ~ 20,000 円 / 募集期間 2 日、取引期間 0 日"

# --- 4. Refresh the "取得日時" timestamp on every data row (2-18). ---
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $newTimestamp
}

# --- 5. Column D got wider. ColumnWidth is in "characters"; the engine
#        adds ~5/6 of a character when it serialises to the raw OOXML
#        <col width>, so back that constant out to land exactly on 41. ---
$ws.Columns.Item(4).ColumnWidth = 41 - 5/6

# --- 6. Rebuild the hyperlinks on F2:F18 so each one targets the right
#        row after the insertion (the COM row-shift above does not slide
#        the worksheet's Hyperlinks collection along with the cells). ---
$ws.Hyperlinks.Delete()

$urls = @{
    2  = "https://www.lancers.jp/work/detail/5422740"
    3  = "https://www.lancers.jp/work/detail/5423046"
    4  = "https://www.lancers.jp/work/detail/5422760"
    5  = "https://www.lancers.jp/work/detail/5422386"
    6  = "https://www.lancers.jp/work/detail/5422331"
    7  = "https://www.lancers.jp/work/detail/5251319"
    8  = "https://www.lancers.jp/work/detail/5422314"
    9  = "https://www.lancers.jp/work/detail/5422652"
    10 = "https://www.lancers.jp/work/detail/5422936"
    11 = "https://www.lancers.jp/work/detail/5422908"
    12 = "https://www.lancers.jp/work/detail/5415980"
    13 = "https://www.lancers.jp/work/detail/5417644"
    14 = "https://www.lancers.jp/work/detail/5423114"
    15 = "https://www.lancers.jp/work/detail/5422916"
    16 = "https://www.lancers.jp/work/detail/5422660"
    17 = "https://www.lancers.jp/work/detail/5420233"
    18 = "https://www.lancers.jp/work/detail/5421782"
}

for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $urls[$r])
    $cell.Style = "Hyperlink"
}
